$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 passive force value tweaks
$ws.Range("B2").Value = 341.91203969887903
$ws.Range("C2").Value = 729.143780703071
$ws.Range("D2").Value = 351.44220921357737
$ws.Range("E2").Value = 553.780641310654

# Row 3 passive force value tweaks
$ws.Range("B3").Value = 405.85417617272054
$ws.Range("C3").Value = 401.84472039471194
$ws.Range("D3").Value = 333.83415277465247
$ws.Range("E3").Value = 345.12933316689646

# Update the selected range to reflect the narrower data of interest
$ws.Range("B1:E3").Select()
